# sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio
$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn sale phụ": chiet khau sale phu giam tu 26000 xuong 16000 ---
$wsSalePhu = $wb.Worksheets.Item("Đơn sale phụ")
$wsSalePhu.Range("N2").Value = 16000
$wsSalePhu.Range("N3").Value = 16000

# --- Sheet "Đơn phụ phẫu 1": them 1 dong don moi (ma 578) truoc dong Tong ---
$wsPhuPhau = $wb.Worksheets.Item("Đơn phụ phẫu 1")
$wsPhuPhau.Rows.Item(9).Insert()

$wsPhuPhau.Range("A9").Value = "HD-LUXURY"
$wsPhuPhau.Range("B9").Value = 578
$wsPhuPhau.Range("C9").NumberFormat = "@"
$wsPhuPhau.Range("C9").Value = "07-21-2024"
$wsPhuPhau.Range("C9").ClearFormats()
$wsPhuPhau.Range("D9").Value = "SÓC TRĂNG"
$wsPhuPhau.Range("E9").Value = "đường thị út"
$wsPhuPhau.Range("F9").Value = "Cá nhân"
$wsPhuPhau.Range("G9").Value = "Nâng cung chân mày"
$wsPhuPhau.Range("H9").Value = "Kha Như Huỳnh "
$wsPhuPhau.Range("I9").Value = 50000

# dong "Tong" (bi day xuong dong 10) cap nhat lai so luong don va tong tien
$wsPhuPhau.Range("B10").Value = 8
$wsPhuPhau.Range("I10").Value = 500000

# --- Sheet "Lương": cap nhat lai cac so lieu tong hop lien quan ---
$wsLuong = $wb.Worksheets.Item("Lương")
$wsLuong.Range("B24").Value = 16000
$wsLuong.Range("B27").Value = 500000
$wsLuong.Range("B33").Value = 4319264.285714285
$wsLuong.Range("B34").Value = 4319264.285714285
